# Updates cryptos list price/volume figures and re-ranks rows 23-51
# (new coin "WrappedliquidstakedEther2.0" inserted at rank 23, shifting
# subsequent coins down by one row and dropping the former last row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''26.913.87'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.46%  '
$ws.Range("D3").Value = '''1.829.16'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.44%  '
$ws.Range("E4").Value = '  +0.65%  '
$ws.Range("D5").Value = '''310.86'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.97%  '
$ws.Range("E6").Value = '  +0.60%  '
$ws.Range("D7").Value = '''0.4580'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.70%  '
$ws.Range("D8").Value = '''0.3691'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.31%  '
$ws.Range("D9").Value = '''0.07182'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.00%  '
$ws.Range("D10").Value = '''0.8767'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.52%  '
$ws.Range("D11").Value = '''0.07845'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.22%  '
$ws.Range("D12").Value = '''19.63'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = '''1.812.32'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.68%  '
$ws.Range("D14").Value = '''5.335'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.86%  '
$ws.Range("D15").Value = '''6.398'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.39%  '
$ws.Range("D16").Value = '''87.19'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.12%  '
$ws.Range("E17").Value = '  +0.63%  '
$ws.Range("D18").Value = '''0.000008718'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.47%  '
$ws.Range("E19").Value = '  +0.51%  '
$ws.Range("D20").Value = '''26.938.77'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.49%  '
$ws.Range("D21").Value = '''14.49'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.16%  '
$ws.Range("D22").Value = '''5.002'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.29%  '
$ws.Range("B23").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D23").Value = '''2.048.34'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.27%  '
$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").Value = '''10.44'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.73%  '
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").Value = '''1.977'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.64%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = '''151.13'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.90%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '''18.20'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.02%  '
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").Value = '''1.968'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.93%  '
$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").Value = '''113.97'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.88%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").Value = '''4.927'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.90%  '
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").Value = '''0.08798'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.66%  '
$ws.Range("B32").Value = 'HuobiToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D32").Value = '''3.052'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.06%  '
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").Value = '''0.7550'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.46%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").Value = '''4.482'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.09%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '''1.132'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.08%  '
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D36").Value = '''2.555'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.37%  '
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").Value = '''1.086'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.30%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '''0.01935'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.07%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '''0.05139'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.32%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = '''2.903'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.53%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = '''6.942'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.66%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '''0.4979'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.44%  '
$ws.Range("B43").Value = 'Algorand'
$ws.Range("C43").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D43").Value = '''0.1597'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.52%  '
$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").Value = '''8.315'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.27%  '
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value = '''0.4684'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.11%  '
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").Value = '''1.007'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.67%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").Value = '''10.18'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.15%  '
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").Value = '''102.04'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.32%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = '''1.616'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.05%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '''0.06116'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.88%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").Value = '''64.35'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.24%  '
